$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
